$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.73"
$ws.Range("E2").Value = "'-4.20%"
$ws.Range("D3").Value = "'30.83"
$ws.Range("E3").Value = "'-6.07%"
$ws.Range("D4").Value = "'4.939"
$ws.Range("E4").Value = "'-0.16%"
$ws.Range("D5").Value = "'0.07202"
$ws.Range("E5").Value = "'-8.22%"
$ws.Range("D6").Value = "'1.787"
$ws.Range("E6").Value = "'-11.94%"
$ws.Range("D7").Value = "'7.667"
$ws.Range("E7").Value = "'-2.17%"
$ws.Range("D8").Value = "'3.748"
$ws.Range("E8").Value = "'-1.59%"
$ws.Range("D9").Value = "'0.8953"
$ws.Range("E9").Value = "'-3.09%"
$ws.Range("D10").Value = "'0.1651"
$ws.Range("E10").Value = "'-5.91%"
$ws.Range("E11").Value = "'-1.45%"
$ws.Range("D12").Value = "'0.08028"
$ws.Range("E12").Value = "'-7.44%"
$ws.Range("D13").Value = "'0.03068"
$ws.Range("E13").Value = "'-2.50%"
$ws.Range("D14").Value = "'0.1002"
$ws.Range("E14").Value = "'-0.29%"
$ws.Range("D15").Value = "'0.001505"
$ws.Range("E15").Value = "'-0.81%"
$ws.Range("D16").Value = "'0.005688"
$ws.Range("E16").Value = "'-3.19%"
$ws.Range("D17").Value = "'3.475"
$ws.Range("E17").Value = "'0.31%"
$ws.Range("D18").Value = "'2.083"
$ws.Range("E18").Value = "'-3.31%"
$ws.Range("D19").Value = "'0.3280"
$ws.Range("E19").Value = "'-0.85%"
$ws.Range("E20").Value = "'-3.35%"
$ws.Range("E21").Value = "'-6.00%"
$ws.Range("D22").Value = "'0.2100"
$ws.Range("E22").Value = "'5.49%"
$ws.Range("D23").Value = "'0.04526"
$ws.Range("E23").Value = "'-0.77%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'-0.65%"
$ws.Range("D25").Value = "'0.004009"
$ws.Range("E25").Value = "'-9.79%"
$ws.Range("E26").Value = "'-0.03%"
$ws.Range("D39").Value = "'0.01599"
$ws.Range("E39").Value = "'-7.97%"
$ws.Range("D40").Value = "'0.04382"
$ws.Range("E40").Value = "'-8.59%"
$ws.Range("D41").Value = "'0.007356"
$ws.Range("E41").Value = "'-1.35%"
$ws.Range("D42").Value = "'0.1307"
$ws.Range("E42").Value = "'-3.92%"
$ws.Range("D43").Value = "'0.007671"
$ws.Range("D44").Value = "'0.002070"
$ws.Range("E44").Value = "'-11.57%"
$ws.Range("D45").Value = "'0.009209"
$ws.Range("E45").Value = "'-12.94%"
$ws.Range("D46").Value = "'0.00005925"
$ws.Range("E46").Value = "'-5.34%"
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("D48").Value = "'2.246"
$ws.Range("E48").Value = "'173.66%"
$ws.Range("E49").Value = "'-3.27%"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E51").Value = "'-0.03%"
